# New week's data row was inserted into the sheet right after the most
# recent entry (row 83), pushing every subsequent row (old 84..148) down
# by one (to 85..149). We reproduce that by inserting a blank row at 84
# (Excel shifts 84..148 -> 85..149 automatically, carrying values/styles
# with them) and then filling the newly-inserted row 84 with the new
# record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(84).Insert()

$ws.Cells.Item(84, 1).Value  = 11
$ws.Cells.Item(84, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(84, 3).Value  = "Bíobío"
$ws.Cells.Item(84, 4).Value  = 45240
$ws.Cells.Item(84, 5).Value  = 8
$ws.Cells.Item(84, 6).Value  = 100112037
$ws.Cells.Item(84, 7).Value  = "Cebollín"
$ws.Cells.Item(84, 8).Value  = "Sin especificar"
$ws.Cells.Item(84, 9).Value  = "Primera"
$ws.Cells.Item(84, 10).Value = 50
$ws.Cells.Item(84, 11).Value = 4000
$ws.Cells.Item(84, 12).Value = 4000
$ws.Cells.Item(84, 13).Value = 4000
$ws.Cells.Item(84, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(84, 15).Value = "Región Metropolitana"
$ws.Cells.Item(84, 16).Value = 111
$ws.Cells.Item(84, 17).Value = 36
$ws.Cells.Item(84, 18).Value = "Hortaliza"
